$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.9238883333333332
$ws.Cells.Item(2, 8).Value = 2.771665
$ws.Cells.Item(2, 9).Value = 0.006623569135547428
$ws.Cells.Item(2, 10).Value = 0.006623569135547428
$ws.Cells.Item(2, 13).Value = 3.5258
$ws.Cells.Item(2, 14).Value = 10.5774
$ws.Cells.Item(2, 15).Value = 0.2648936965472837
$ws.Cells.Item(2, 16).Value = 0.2648936965472837
$ws.Cells.Item(2, 17).Value = 3.257445485666667
$ws.Cells.Item(2, 18).Value = 29.317009371
$ws.Cells.Item(2, 19).Value = 0.001754541712651655
$ws.Cells.Item(2, 20).Value = 0.001754541712651655

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.9238883333333332
$ws.Cells.Item(3, 8).Value = 2.771665
$ws.Cells.Item(3, 9).Value = 0.006623569135547428
$ws.Cells.Item(3, 10).Value = 0.006623569135547428
$ws.Cells.Item(3, 15).Value = 0.4949431058556773
$ws.Cells.Item(3, 16).Value = 0.4949431058556773
$ws.Cells.Item(3, 17).Value = 6.086404496769998
$ws.Cells.Item(3, 18).Value = 54.77764047092999
$ws.Cells.Item(3, 19).Value = 0.003278289879797647
$ws.Cells.Item(3, 20).Value = 0.003278289879797647

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.9238883333333332
$ws.Cells.Item(4, 8).Value = 2.771665
$ws.Cells.Item(4, 9).Value = 0.006623569135547428
$ws.Cells.Item(4, 10).Value = 0.006623569135547428
$ws.Cells.Item(4, 15).Value = 0.240163197597039
$ws.Cells.Item(4, 16).Value = 0.240163197597039
$ws.Cells.Item(4, 17).Value = 2.953330086871666
$ws.Cells.Item(4, 18).Value = 26.579970781845
$ws.Cells.Item(4, 19).Value = 0.001590737543098126
$ws.Cells.Item(4, 20).Value = 0.001590737543098126

# Row 5
$ws.Cells.Item(5, 9).Value = 0.9911582946048672
$ws.Cells.Item(5, 10).Value = 0.9911582946048673
$ws.Cells.Item(5, 13).Value = 3.5258
$ws.Cells.Item(5, 14).Value = 10.5774
$ws.Cells.Item(5, 15).Value = 0.2648936965472837
$ws.Cells.Item(5, 16).Value = 0.2648936965472837
$ws.Cells.Item(5, 17).Value = 487.4477862719333
$ws.Cells.Item(5, 18).Value = 4387.0300764474
$ws.Cells.Item(5, 19).Value = 0.2625515845213849
$ws.Cells.Item(5, 20).Value = 0.2625515845213849

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9911582946048672
$ws.Cells.Item(6, 10).Value = 0.9911582946048673
$ws.Cells.Item(6, 15).Value = 0.4949431058556773
$ws.Cells.Item(6, 16).Value = 0.4949431058556773
$ws.Cells.Item(6, 19).Value = 0.4905669647263494
$ws.Cells.Item(6, 20).Value = 0.4905669647263494

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9911582946048672
$ws.Cells.Item(7, 10).Value = 0.9911582946048673
$ws.Cells.Item(7, 15).Value = 0.240163197597039
$ws.Cells.Item(7, 16).Value = 0.240163197597039
$ws.Cells.Item(7, 19).Value = 0.2380397453571329
$ws.Cells.Item(7, 20).Value = 0.2380397453571329

# Row 8
$ws.Cells.Item(8, 8).Value = 0.92819
$ws.Cells.Item(8, 9).Value = 0.0022181362595854
$ws.Cells.Item(8, 10).Value = 0.0022181362595854
$ws.Cells.Item(8, 13).Value = 3.5258
$ws.Cells.Item(8, 14).Value = 10.5774
$ws.Cells.Item(8, 15).Value = 0.2648936965472837
$ws.Cells.Item(8, 16).Value = 0.2648936965472837
$ws.Cells.Item(8, 17).Value = 1.090870767333333
$ws.Cells.Item(8, 18).Value = 9.817836906
$ws.Cells.Item(8, 19).Value = 0.0005875703132471419
$ws.Cells.Item(8, 20).Value = 0.0005875703132471419

# Row 9
$ws.Cells.Item(9, 8).Value = 0.92819
$ws.Cells.Item(9, 9).Value = 0.0022181362595854
$ws.Cells.Item(9, 10).Value = 0.0022181362595854
$ws.Cells.Item(9, 15).Value = 0.4949431058556773
$ws.Cells.Item(9, 16).Value = 0.4949431058556773
$ws.Cells.Item(9, 19).Value = 0.001097851249530293
$ws.Cells.Item(9, 20).Value = 0.001097851249530293

# Row 10
$ws.Cells.Item(10, 8).Value = 0.92819
$ws.Cells.Item(10, 9).Value = 0.0022181362595854
$ws.Cells.Item(10, 10).Value = 0.0022181362595854
$ws.Cells.Item(10, 15).Value = 0.240163197597039
$ws.Cells.Item(10, 16).Value = 0.240163197597039
$ws.Cells.Item(10, 17).Value = 0.9890269759633333
$ws.Cells.Item(10, 19).Value = 0.0005327146968079654
$ws.Cells.Item(10, 20).Value = 0.0005327146968079654
